$wb = $excel.ActiveWorkbook

# --- Update Metadata sheet timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("A2").Value = "29 Oct 2025, 10:22 AM"

# --- Update Top Gainers sheet rows 65-76 ---
$gainers = $wb.Worksheets.Item("Top Gainers")
$gainers.Range("B65").Value = "ICRA"
$gainers.Range("C65").Value = 3.7985
$gainers.Range("D65").Value = 4.4793
$gainers.Range("E65").Value = 2.8828

$gainers.Range("B66").Value = "SALASAR"
$gainers.Range("C66").Value = 3.7935
$gainers.Range("D66").Value = 4.7872
$gainers.Range("E66").Value = 11.0485

$gainers.Range("B67").Value = "NPST"
$gainers.Range("C67").Value = 3.7841
$gainers.Range("D67").Value = -2.0689
$gainers.Range("E67").Value = -3.5677

$gainers.Range("B68").Value = "DCW"
$gainers.Range("C68").Value = 3.7544
$gainers.Range("D68").Value = 2.3219
$gainers.Range("E68").Value = -3.9753

$gainers.Range("B69").Value = "RHETAN"
$gainers.Range("C69").Value = 3.754
$gainers.Range("D69").Value = 4.178
$gainers.Range("E69").Value = 6.549

$gainers.Range("B70").Value = "HINDPETRO"
$gainers.Range("C70").Value = 3.6935
$gainers.Range("D70").Value = 6.9335
$gainers.Range("E70").Value = 5.7397

$gainers.Range("B71").Value = "SHK"
$gainers.Range("C71").Value = 3.6851
$gainers.Range("D71").Value = 2.4377
$gainers.Range("E71").Value = -1.8843

$gainers.Range("B72").Value = "BHARTIHEXA"
$gainers.Range("C72").Value = 3.6718
$gainers.Range("D72").Value = 7.0877
$gainers.Range("E72").Value = 15.3332

$gainers.Range("B73").Value = "HLEGLAS"
$gainers.Range("C73").Value = 3.659
$gainers.Range("D73").Value = 8.115500000000001
$gainers.Range("E73").Value = 27.1239

$gainers.Range("B74").Value = "RHIM"
$gainers.Range("C74").Value = 3.6544
$gainers.Range("D74").Value = 3.2276
$gainers.Range("E74").Value = 5.1826

$gainers.Range("B75").Value = "BCLIND"
$gainers.Range("C75").Value = 3.6271
$gainers.Range("D75").Value = 2.2945
$gainers.Range("E75").Value = 0.1728

$gainers.Range("B76").Value = "CGPOWER"
$gainers.Range("C76").Value = 3.6125
$gainers.Range("D76").Value = 3.4192
$gainers.Range("E76").Value = 1.0325
